$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.12579619884491
$ws.Range("B1").Value = 1.118985891342163
$ws.Range("C1").Value = 5.23612117767334
$ws.Range("D1").Value = 2.122360706329346
$ws.Range("E1").Value = 1.201128840446472
